$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.597641587257385
$ws.Range("B1").Value = 1.925802946090698
$ws.Range("C1").Value = 2.038169622421265
$ws.Range("D1").Value = 2.395359754562378
$ws.Range("E1").Value = 3.281719446182251
